# Updated cryptos list on Thu Jun  6 10:14:59 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $coin, $link, $price, $volume) {
    if ($coin -ne $null) {
        $ws.Cells.Item($row, 2).Value = $coin
    }
    if ($link -ne $null) {
        $ws.Cells.Item($row, 3).Value = $link
    }
    if ($price -ne $null) {
        $ws.Cells.Item($row, 4).Value = $price
    }
    if ($volume -ne $null) {
        $ws.Cells.Item($row, 5).Value = $volume
    }
}

Set-Row 2  $null $null "70.874.45"   "  +0.02%  "
Set-Row 3  $null $null "3.846.01"    "  +1.32%  "
Set-Row 4  $null $null $null          "  -0.01%  "
Set-Row 5  $null $null "705.92"      "  +1.16%  "
Set-Row 6  $null $null "172.42"      "  -0.57%  "
Set-Row 7  $null $null "3.843.65"    "  +1.31%  "
Set-Row 8  $null $null $null          "  -0.02%  "
Set-Row 9  $null $null "0.523"       "  -0.67%  "
Set-Row 10 $null $null $null          "  -0.48%  "
Set-Row 11 $null $null "7.28"        "  -2.29%  "
Set-Row 12 $null $null "0.457"       "  -0.54%  "
Set-Row 13 $null $null $null          "  -0.46%  "
Set-Row 14 $null $null "36.58"       "  +0.73%  "
Set-Row 15 $null $null "4.493.92"    "  +1.32%  "
Set-Row 16 $null $null "3.891.80"    "  +2.43%  "
Set-Row 17 $null $null "70.913.98"   "  +0.05%  "
Set-Row 18 $null $null "7.19"        "  -0.15%  "
Set-Row 19 $null $null $null          "  +0.63%  "
Set-Row 20 $null $null "17.36"       "  -2.63%  "
Set-Row 21 $null $null "491.48"      "  +2.02%  "
Set-Row 22 $null $null "10.60"       "  -4.89%  "
Set-Row 23 $null $null "0.716"       "  +0.32%  "
Set-Row 24 $null $null "85.07"       "  +0.91%  "
Set-Row 25 $null $null $null          "  +2.27%  "
Set-Row 26 $null $null "10.55"       "  +1.10%  "
Set-Row 27 $null $null "12.13"       "  -2.29%  "
Set-Row 28 $null $null $null          "  -2.42%  "
Set-Row 29 $null $null $null          "  +3.27%  "
Set-Row 30 $null $null $null          "  -0.04%  "
Set-Row 31 $null $null $null          "  -0.30%  "
Set-Row 32 $null $null $null          "  -0.39%  "

# Rows 33/34 swap: Kaspa <-> EthereumClassic
Set-Row 33 "EthereumClassic" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc" "29.36" "  -0.79%  "
Set-Row 34 "Kaspa" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas" "0.180" "  -0.32%  "

Set-Row 35 $null $null "9.17"        "  -0.50%  "
Set-Row 36 $null $null "3.801.72"    "  +1.50%  "
Set-Row 37 $null $null $null          "  +0.49%  "
Set-Row 38 $null $null "0.103"       "  +0.29%  "
Set-Row 39 $null $null $null          "  +7.11%  "
Set-Row 40 $null $null "6.05"        "  +0.66%  "
Set-Row 41 $null $null $null          "  +6.43%  "
Set-Row 42 $null $null "3.32"        "  -4.74%  "
Set-Row 43 $null $null $null          "  +0.00%  "
Set-Row 44 $null $null $null          "  +0.16%  "
Set-Row 45 $null $null $null          "  -5.16%  "
Set-Row 46 $null $null "163.36"      "  +0.38%  "
Set-Row 47 $null $null "48.66"       "  -0.73%  "

# Rows 48/49 swap: ONDO <-> Bittensor
Set-Row 48 "Bittensor" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao" "418.99" "  +3.67%  "
Set-Row 49 "ONDO" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo" "1.39" "  +0.53%  "

Set-Row 50 $null $null "0.298"       "  -0.89%  "
Set-Row 51 $null $null "8.61"        "  +0.65%  "
